$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16 and 17 swap order (ShibaInu moves above WrappedEther) and get new price/volume data
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000137"
$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.674.86"
$ws.Range("E17").Value = "  +0.05%  "

# Updated prices and 1h volume percentages for remaining rows
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.505.13"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.650.94"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.83"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.21"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.573"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.71"
$ws.Range("E9").Value = "  -3.62%  "
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.339"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.119.16"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.476.26"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.94"
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "342.76"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.46"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.64"
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.41"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.80"
$ws.Range("E23").Value = "  +2.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.419"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.38"
$ws.Range("E30").Value = "  -4.49%  "
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.07"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "152.05"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.18"
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.862"
$ws.Range("E36").Value = "  -4.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.864"
$ws.Range("E37").Value = "  -4.15%  "
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.57"
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0979"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  -3.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "271.60"
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.47"
$ws.Range("E45").Value = "  -3.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.71"
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0537"
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.045.97"
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.78"
$ws.Range("E49").Value = "  -3.07%  "
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.81"
$ws.Range("E51").Value = "  -2.30%  "

Write-Host "Cryptos list updated"
